$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DM_Stat (column C) and P_Value (column D) values per row
$ws.Range("C2").Value = 0.8409380881357338
$ws.Range("D2").Value = 0.409430119213916

$ws.Range("C3").Value = 0.3646170402033728
$ws.Range("D3").Value = 0.7188787084779382

$ws.Range("C4").Value = 0.1199620767861547
$ws.Range("D4").Value = 0.9056019817740011

$ws.Range("C5").Value = 0.04129806189883264
$ws.Range("D5").Value = 0.9674307903600394

$ws.Range("C6").Value = -0.7023825779698128
$ws.Range("D6").Value = 0.4898062167577877

$ws.Range("C7").Value = -1.023094058607468
$ws.Range("D7").Value = 0.3173824246139023

$ws.Range("C8").Value = -0.8889404025704247
$ws.Range("D8").Value = 0.3836521469958374

$ws.Range("C9").Value = -0.2165376456570488
$ws.Range("D9").Value = 0.830564015577361

$ws.Range("C10").Value = -0.3401071768065458
$ws.Range("D10").Value = 0.7370017346845437

$ws.Range("C11").Value = -0.04333921033373841
$ws.Range("D11").Value = 0.9658220922184564
